$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Astronauta" (sheet1): add attendance marks in column D
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Astronauta")
$ws1.Activate()

$ws1.Range("D2").Value = 1

# D3 must be stored as literal TEXT "0.75" (so it is excluded from the
# COUNT()/SUM() used by column L), not as a real number.
$ws1.Range("D3").NumberFormat = "@"
$ws1.Range("D3").Value = "0.75"
$ws1.Range("D3").Style = "Normal"

$ws1.Range("D7").Value = 1

[void]$ws1.Range("D8").Select()

# ---------------------------------------------------------------------------
# Sheet "Senador" (sheet2): add a single new attendance mark in column D
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Senador")
$ws2.Activate()

$ws2.Range("D4").Value = 0

[void]$ws2.Range("A2:A7").Select()

# ---------------------------------------------------------------------------
# Sheet "Mago" (sheet3): add attendance marks in columns D, E and F
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Mago")
$ws3.Activate()

$ws3.Range("D2").Value = 1
$ws3.Range("E2").Value = 1

$ws3.Range("D3").Value = 1

$ws3.Range("D4").Value = 1
$ws3.Range("E4").Value = 0

$ws3.Range("E5").Value = 1

$ws3.Range("D6").Value = 1
$ws3.Range("E6").Value = 1
$ws3.Range("F6").Value = 1

$ws3.Range("D7").Value = 1

[void]$ws3.Range("G5").Select()

# ---------------------------------------------------------------------------
# Sheet "Ninja" (sheet4): add attendance marks in column E
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Ninja")
$ws4.Activate()

$ws4.Range("E2").Value = 1
$ws4.Range("E3").Value = 0

# E4 must be stored as literal TEXT "0.5", excluded from the P column average.
$ws4.Range("E4").NumberFormat = "@"
$ws4.Range("E4").Value = "0.5"
$ws4.Range("E4").Style = "Normal"

$ws4.Range("E5").Value = 0
$ws4.Range("E7").Value = 0

[void]$ws4.Range("F4").Select()
